$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "marker_2" column header
$ws.Range("K1").Value = "marker_2"

# Enter marker info for the two off-by-one rows (NAT marker)
$ws.Range("K15").Value = "NAT"
$ws.Range("K26").Value = "NAT"

# Move active selection to E15, matching the saved view state
$ws.Range("E15").Select()
